# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (volume number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Cells that change data type (text <-> numeric): copy style/value from a donor cell first ---
$ws.Range("J14").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("J14").Copy($ws.Range("F14"))
$ws.Range("F14").Value = 1
$ws.Range("J14").Copy($ws.Range("I14"))
$ws.Range("I14").Value = 1
$ws.Range("D14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("D14").Copy($ws.Range("C18"))
$ws.Range("J14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("J14").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("K16").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100
$ws.Range("J14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 2
$ws.Range("K16").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("D14").Copy($ws.Range("G31"))
$ws.Range("E14").Copy($ws.Range("H31"))

# --- Plain value updates (same type/style, only the number changes) ---
$ws.Range("K14").Value = -50
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -66.666666666666
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = 112.5
$ws.Range("L15").Value = 112.5
$ws.Range("M15").Value = -5.555555555555
$ws.Range("N15").Value = -5.555555555555
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = -50
$ws.Range("L16").Value = -16.39344262295
$ws.Range("M16").Value = -53.636363636363
$ws.Range("N16").Value = -82.474226804123
$ws.Range("F17").Value = 18
$ws.Range("H17").Value = 80
$ws.Range("I17").Value = 158
$ws.Range("J17").Value = 159
$ws.Range("K17").Value = -0.62893081761
$ws.Range("L17").Value = 30.578512396694
$ws.Range("M17").Value = 23.4375
$ws.Range("N17").Value = -45.328719723183
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -81.25
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = 1.075268817204
$ws.Range("L18").Value = 22.077922077922
$ws.Range("M18").Value = -52.525252525252
$ws.Range("N18").Value = -92.610062893081
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 166.666666666667
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 437
$ws.Range("J19").Value = 446
$ws.Range("K19").Value = -2.01793721973
$ws.Range("L19").Value = 25.936599423631
$ws.Range("M19").Value = 12.919896640826
$ws.Range("N19").Value = -44.892812105926
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 73
$ws.Range("J20").Value = 105
$ws.Range("K20").Value = -30.47619047619
$ws.Range("L20").Value = -43.846153846153
$ws.Range("M20").Value = -34.234234234234
$ws.Range("N20").Value = -97.232752084912
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 7.692307692307
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = -12.698412698412
$ws.Range("I21").Value = 831
$ws.Range("J21").Value = 860
$ws.Range("K21").Value = -3.372093023255
$ws.Range("L21").Value = 11.543624161073
$ws.Range("M21").Value = -12.893081761006
$ws.Range("N21").Value = -84.33257918552
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 33
$ws.Range("K23").Value = -21.212121212121
$ws.Range("L23").Value = 30
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -44.444444444444
$ws.Range("F24").Value = 56
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = -39.130434782608
$ws.Range("I24").Value = 936
$ws.Range("J24").Value = 1035
$ws.Range("K24").Value = -9.565217391304
$ws.Range("L24").Value = 15.4130702836
$ws.Range("M24").Value = -41.242937853107
$ws.Range("C25").Value = 11
$ws.Range("E25").Value = 83.333333333333
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 461
$ws.Range("J25").Value = 412
$ws.Range("K25").Value = 11.893203883495
$ws.Range("L25").Value = 118.483412322275
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -35.483870967741
$ws.Range("I26").Value = 338
$ws.Range("J26").Value = 308
$ws.Range("K26").Value = 9.740259740259
$ws.Range("L26").Value = 6.2893081761
$ws.Range("M26").Value = -33.725490196078
$ws.Range("F27").Value = 1
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = 57.142857142857
$ws.Range("L27").Value = 15.78947368421
$ws.Range("G28").Value = 4
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = -10.81081081081
$ws.Range("L33").Value = 100
